# fix: deidentify personal emails
# Replace the "roland.stens..." test fixture emails with "test.user..." ones
# across the "Teams" sheet, and keep the view scrolled to column M (matches
# the author's saved sheetView topLeftCell="M1").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Teams")

$dim = $ws.Range("A1:Z8")

for ($r = 1; $r -le $dim.Rows.Count; $r++) {
    for ($c = 1; $c -le $dim.Columns.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $val -is [string] -and $val.Contains("roland.stens")) {
            $cell.Value = $val.Replace("roland.stens", "test.user")
        }
    }
}

# Match the author's saved scroll position on the Teams sheet view.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 13
